$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 (subject numbers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON)
$ws.Range("B2").Value = 461.62263204731903
$ws.Range("C2").Value = 389.61793576387959
$ws.Range("D2").Value = 457.90149340365843
$ws.Range("E2").Value = 385.30257189742673

# Row 3 (STR)
$ws.Range("B3").Value = 463.12862074754423
$ws.Range("C3").Value = 388.3339396377682
$ws.Range("D3").Value = 458.30069618978428
$ws.Range("E3").Value = 395.4634356969558

# Update selection to match new edited range
$ws.Range("B1:E3").Select()
